$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting (borders, alignment, etc.) from the row above into the new row
$ws.Range("A10:C10").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add new row of data (row 11) mirroring the style/format of the previous sprint rows
$ws.Range("A11").Value = 44643
$ws.Range("B11").Value = "Le septième sprint a comencé"
$ws.Range("C11").Value = ""

# Resize the table to include the new row
$table = $ws.ListObjects.Item("Tableau2")
$table.Resize($ws.Range("A1:C11"))

# Update the active selection cell to mirror the new state
$ws.Range("B19").Select()

$wb.Save()
